$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Sheet1: move the selection to B2 (it will no longer be the active tab
# once Sheet2 is added/activated below).
$null = $sheet1.Activate()
$null = $sheet1.Range("B2").Select()

# Insert the new "Sheet2" right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# Header row.
$ws2.Range("A1").Value = "Input"
$ws2.Range("B1").Value = "Equ"
$ws2.Range("C1").Value = "Equ``"

# Newton's Method example data + formulas.
$ws2.Range("A2").Value = 10
$ws2.Range("B2").Formula = "=A2^6 -A2^5-6*A2^4-A2^2+A2+10"
$ws2.Range("C2").Formula = "=6*A2^5-5*A2^4-24*A2^3-2*A2+1"
$ws2.Range("A3").Formula = "=A2-B2/C2"

# Column widths (bestFit-style) for columns B and C.
$ws2.Columns.Item(2).ColumnWidth = 11.16
$ws2.Columns.Item(3).ColumnWidth = 10.16

# Sheet2 becomes the active sheet/tab with this selection.
$null = $ws2.Activate()
$null = $ws2.Range("C9").Select()
